# Impediment_Log_Sprint 3.xlsx — Sprint 3 impediment log update
# - IMP-007: tighten the resolution text and correct the Raised On date
# - IMP-008: mark as closed and stamp the Resolved On date
# - IMP-009: brand-new impediment row (architecture change, DB -> cloud)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- IMP-007 (row 21): shorten resolution text, fix "Raised On" date ---
$ws.Range("E21").Value = 45970
$ws.Range("I21").Value = "Revisar el método de encuesta"

# --- IMP-008 (row 22): now resolved -> Cerrado, with a Resolved On date ---
$ws.Range("H22").Value = "Cerrado"
$ws.Range("J21").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 45988

# --- IMP-009 (new row 23): architecture change from DB to cloud storage ---
$ws.Rows(23).Insert()
$ws.Range("B22:K22").Copy()
$ws.Range("B23:K23").PasteSpecial(-4122)

$ws.Range("B23").Value = "IMP-009"
$ws.Range("C23").Value = "Modificacion de Arquitectura de Capa Backend"
$ws.Range("D23").Value = "Tester"
$ws.Range("E23").Value = 45958
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = "Dev"
$ws.Range("H23").Value = "Cerrado"
$ws.Range("I23").Value = "Se cambió la Base de Datos Tradicional por un Sistema de almacenamiento en Nube y controlador Lógico de Nube"
$ws.Range("I23").WrapText = $true
$ws.Range("J23").Value = 45967
$ws.Rows(23).RowHeight = 27

# --- view state: scroll position + current selection ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 12
$win.ScrollColumn = 6
$ws.Range("I28").Select()
